# Add the two new mnemonic entries (X_COORD / Y_COORD) used for the
# municipality-seat coordinates coming from the shapefile, following the
# existing "mnemonico -> descricao" layout on the gas_prices_hist sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write the short codes first, then the long descriptions, so the shared
# string table is populated in the same order as the source workbook
# (X_COORD, Y_COORD, then their two descriptions).
$ws.Range("A42").Value = "X_COORD"
$ws.Range("A43").Value = "Y_COORD"
$ws.Range("B42").Value = "X cordenada da sede do município"
$ws.Range("B43").Value = "Y cordenada da sede do município"

# Reflect the cursor ending up on the row right after the new entries,
# matching the saved view state of the edited workbook.
$ws.Range("B44").Select()
